$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F updates: row -> new value
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 9944
    3  = 419
    6  = 278
    9  = 738
    11 = 1222
    12 = 1013
    13 = 3082
    14 = 2324
    16 = 2024
    19 = 485
    20 = 1575
    21 = 533
    23 = 213
    27 = 74
    28 = 347
    29 = 555
    30 = 44
    31 = 205
    32 = 1562
    33 = 273
    34 = 1612
    35 = 87
    36 = 391
    38 = 422
    39 = 888
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# Sheet "演出" (shows) - column F updates: row -> new value
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    4 = 32
    9 = 5
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item($row, 6).Value = $showUpdates[$row]
}

# Sheet "全部类型" (all types, combined) - column F updates: row -> new value
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 9944
    3  = 419
    8  = 278
    11 = 738
    13 = 1222
    14 = 1013
    15 = 3082
    16 = 2324
    18 = 2024
    21 = 485
    22 = 1575
    23 = 533
    25 = 213
    29 = 74
    30 = 347
    31 = 555
    32 = 32
    35 = 44
    36 = 205
    37 = 1562
    39 = 274
    40 = 1612
    41 = 87
    43 = 391
    45 = 422
    46 = 888
    49 = 5
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
